# Apply updated "想去人数" (want-to-go count) values across sheets.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 9445
$ws1.Range("F7").Value = 204
$ws1.Range("F21").Value = 369
$ws1.Range("F26").Value = 273
$ws1.Range("F31").Value = 630
$ws1.Range("F36").Value = 312
$ws1.Range("F37").Value = 510
$ws1.Range("F40").Value = 510
$ws1.Range("F45").Value = 312

# Sheet "演出" (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 16
$ws2.Range("F24").Value = 7
$ws2.Range("F35").Value = 112

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 9445
$ws4.Range("F22").Value = 369
$ws4.Range("F30").Value = 273
$ws4.Range("F36").Value = 630
$ws4.Range("F48").Value = 112
